$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 322 (shifts existing rows 322.. down by one)
$ws.Rows(322).Insert()

# Fill in the constant columns (same for every data row in this sheet)
$ws.Cells.Item(322, 1).Value = 9
$ws.Cells.Item(322, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(322, 3).Value = "Metropolitana"
$ws.Cells.Item(322, 4).Value = 44726
$ws.Cells.Item(322, 5).Value = 13
$ws.Cells.Item(322, 6).Value = 100112012
$ws.Cells.Item(322, 7).Value = "Espinaca"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 210
$ws.Cells.Item(322, 11).Value = 7500
$ws.Cells.Item(322, 12).Value = 8000
$ws.Cells.Item(322, 13).Value = 7786
$ws.Cells.Item(322, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(322, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(322, 16).Value = 779
$ws.Cells.Item(322, 17).Value = 10
$ws.Cells.Item(322, 18).Value = "Hortaliza"
